# Generate Report for Handback
#
# This applies the "handback" update to the localization-status workbook:
#  - Overview sheet: status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US" (affects all cells sharing that text)
#    and the two status columns (E/F) are widened.
#  - zh-cn / de-de sheets: the "Latest Target File" (I) and
#    "Latest Handback File" (J) columns are populated for both data rows,
#    a hyperlink (pointing at the same source .md file as column A) is
#    added on the "Latest Target File" cell, the "Latest Handback
#    DateTime" (K) is stamped with the handback time, and columns
#    C/I/J are widened to fit the new content.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

# Update status text wherever it currently reads "Ready for handoff".
$statusRanges = @("E2", "F2", "E3", "F3")
foreach ($addr in $statusRanges) {
    $cell = $overview.Range($addr)
    if ($cell.Value() -eq "Ready for handoff") {
        $cell.Value = "Handed back: in sync with en-US"
    }
}

# Widen the zh-cn / de-de status columns to fit the longer text.
$overview.Range("E1").ColumnWidth = 29.2
$overview.Range("F1").ColumnWidth = 29.2

# ---------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de)
# ---------------------------------------------------------------------
$languages = @(
    @{ Sheet = "zh-cn"; Suffix = "zh-cn"; HandbackTime = "2016-08-21 20:58:44" },
    @{ Sheet = "de-de"; Suffix = "de-de"; HandbackTime = "2016-08-21 20:58:51" }
)

$files = @(
    @{ Row = 2; Id = "553ec807-63b4-43f0-bcd3-026bbf09903a"; Hash = "33c48bdf634c03775f97a472f62de3d3b31b1da7" },
    @{ Row = 3; Id = "79f253f7-18fb-4513-98e6-40711ea9adf9"; Hash = "9bfa01776a963d3c12df22fe5d365d2d7ba03488" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Widen the "Status" (C), "Latest Target File" (I) and
    # "Latest Handback File" (J) columns.
    $ws.Range("C1").ColumnWidth = 29.2
    $ws.Range("I1").ColumnWidth = 39.2
    $ws.Range("J1").ColumnWidth = 39.2

    foreach ($file in $files) {
        $row = $file.Row
        $mdName = "$($file.Id).md"
        $targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bf984ccf4a0879c7fc8bee8ee13e1b02cc0ff384/e2e/$mdName"
        $handbackFile = "$($file.Id).$($file.Hash).$($lang.Suffix).xlf"

        $iCell = $ws.Range("I$row")
        $jCell = $ws.Range("J$row")
        $kCell = $ws.Range("K$row")

        # Latest Target File: same file as the source column (A), with a
        # matching hyperlink.
        $iCell.Value = $mdName
        $ws.Hyperlinks.Add($iCell, $targetUrl, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null

        # Latest Handback File.
        $jCell.Value = $handbackFile

        # Latest Handback DateTime.
        $kCell.Value = $lang.HandbackTime
    }
}
